# Swap the data of rows 86 and 87 (the match records for id 7511958 and
# 7511976) while leaving the shared columns (A, C, D, E, J) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 86
$row2 = 87

# Columns that hold per-record data which must be swapped between the two
# rows. (A, C, D, E, J are identical between the two rows and stay as-is.)
$cols = @("B","F","G","H","I","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"

    $val1 = $ws.Range($addr1).Value2
    $val2 = $ws.Range($addr2).Value2

    $ws.Range($addr1).Value2 = $val2
    $ws.Range($addr2).Value2 = $val1
}
